# Dev 20012024 (Post pre-go-live workshop)
#
# The "survey" sheet had two fields removed from the "create_new_person"
# group: the now-unused "external_id"/"External ID" row and the
# "geolocation" calculate row. Deleting both rows shifts every row below
# them up by two (which also accounts for the shrunk used-range / dropped
# trailing blank rows at the bottom of the sheet).
#
# The health_unit "Set the Primary Contact" note's calculate expression
# also had its separator text tweaked from "'s Primary Contact" to
# " - Primary Contact".
#
# Finally, the reviewer's cursor/scroll position when they saved (visible
# in the sheetView pane/selection) moved from C5/C20 to J10/J19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Activate() | Out-Null

# Remove the "external_id" row (old row 32) and, once everything has
# shifted up, the "geolocation" row (old row 34, now row 33).
$ws.Rows("32:32").Delete() | Out-Null
$ws.Rows("33:33").Delete() | Out-Null

# Update the calculate expression on the health_unit primary-contact note.
$ws.Range("J18").Value = 'concat(../../health_unit/name, " - Primary Contact")'

# Reflect the reviewer's final cursor position/scroll in the saved view.
$ws.Range("J19").Select() | Out-Null
